# network_informations.xlsx - "Modifiction du tableau Excel"
#
# 1) "interfaces" sheet: fix the stray "Branch3" label on row 15 (should be
#    "Branch2", matching the other S0/0/1 rows) and fill in the previously
#    empty " Adresse IP" / " Masque de sous-réseau" columns (D:E) for every
#    interface row.
# 2) "static routes" sheet: populate the (until now empty) static routing
#    table with the Routeur / destination network / mask / gateway rows.
# 3) Make "static routes" the active sheet/tab, matching the new selections
#    recorded on each sheet.

$wb = $excel.ActiveWorkbook

$wsNetworks   = $wb.Worksheets.Item("networks")
$wsInterfaces = $wb.Worksheets.Item("interfaces")
$wsRoutes     = $wb.Worksheets.Item("static routes")

# ---------------------------------------------------------------------------
# 1) interfaces
# ---------------------------------------------------------------------------

# Correct mislabeled device on row 15 (was "Branch3", should be "Branch2")
$wsInterfaces.Range("B15").Value = "Branch2"

# Fill in IP address (D) and subnet mask (E) for every interface row 4-15
$wsInterfaces.Range("D4").Value = "192.168.1.1"
$wsInterfaces.Range("E4").Value = "255.255.255.192"

$wsInterfaces.Range("D5").Value = "192.168.1.65"
$wsInterfaces.Range("E5").Value = "255.255.255.192"

$wsInterfaces.Range("D6").Value = "192.168.1.225"
$wsInterfaces.Range("E6").Value = "255.255.255.252"

$wsInterfaces.Range("D7").Value = "192.168.1.229"
$wsInterfaces.Range("E7").Value = "255.255.255.252"

$wsInterfaces.Range("D8").Value = "192.168.1.129"
$wsInterfaces.Range("E8").Value = "255.255.255.224"

$wsInterfaces.Range("D9").Value = "192.168.1.161"
$wsInterfaces.Range("E9").Value = "255.255.255.224"

$wsInterfaces.Range("D10").Value = "192.168.1.226"
$wsInterfaces.Range("E10").Value = "255.255.255.252"

$wsInterfaces.Range("D11").Value = "192.168.1.233"
$wsInterfaces.Range("E11").Value = "255.255.255.252"

$wsInterfaces.Range("D12").Value = "192.168.1.193"
$wsInterfaces.Range("E12").Value = "255.255.255.240"

$wsInterfaces.Range("D13").Value = "192.168.1.209"
$wsInterfaces.Range("E13").Value = "255.255.255.240"

$wsInterfaces.Range("D14").Value = "192.168.1.234"
$wsInterfaces.Range("E14").Value = "255.255.255.252"

$wsInterfaces.Range("D15").Value = "192.168.1.230"
$wsInterfaces.Range("E15").Value = "255.255.255.252"

# ---------------------------------------------------------------------------
# 2) static routes
# ---------------------------------------------------------------------------

$routes = @(
    @("HQ", "192.168.1.160", "255.255.255.224", "192.168.1.226"),
    @("HQ", "192.168.1.128", "255.255.255.224", "192.168.1.226"),
    @("HQ", "192.168.1.208", "255.255.255.240", "192.168.1.230"),
    @("HQ", "192.168.1.192", "255.255.255.240", "192.168.1.30"),
    @("HQ", "192.168.1.232", "255.255.255.252", "192.168.1.226"),
    @("B1", "192.168.1.0",   "255.255.255.192", "192.168.1.225"),
    @("B1", "192.168.1.64",  "255.255.255.192", "192.168.1.225"),
    @("B1", "192.168.1.192", "255.255.255.240", "192.168.1.234"),
    @("B1", "192.168.1.208", "255.255.255.240", "192.168.1.234"),
    @("B1", "192.168.1.228", "255.255.255.252", "192.168.1.234"),
    @("B2", "192.168.1.128", "255.255.255.224", "192.168.1.233"),
    @("B2", "192.168.1.160", "255.255.255.224", "192.168.1.233"),
    @("B2", "192.168.1.0",   "255.255.255.192", "192.168.1.229"),
    @("B2", "192.168.1.64",  "255.255.255.192", "192.168.1.229"),
    @("B2", "192.168.1.224", "255.255.255.252", "192.168.1.233")
)

$row = 4
foreach ($r in $routes) {
    $wsRoutes.Cells.Item($row, 2).Value = $r[0]
    $wsRoutes.Cells.Item($row, 3).Value = $r[1]
    $wsRoutes.Cells.Item($row, 4).Value = $r[2]
    $wsRoutes.Cells.Item($row, 5).Value = $r[3]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3) Selections / active sheet / zoom
# ---------------------------------------------------------------------------

[void]$wsNetworks.Range("B4").Select()

[void]$wsInterfaces.Activate()
[void]$wsInterfaces.Range("D4").Select()
$excel.ActiveWindow.Zoom = 150

[void]$wsRoutes.Activate()
[void]$wsRoutes.Range("F21").Select()
$excel.ActiveWindow.Zoom = 122
